# Refresh the crypto price/volume snapshot (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.992.89"
$ws.Range("E2").Value = "  +2.01%  "

$ws.Range("D3").Value = "'1.643.59"
$ws.Range("E3").Value = "  +0.38%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").Value = "'212.72"
$ws.Range("E5").Value = "  +0.27%  "

$ws.Range("E6").Value = "  -1.15%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.15%  "

$ws.Range("D8").Value = "'23.41"
$ws.Range("E8").Value = "  +2.06%  "

$ws.Range("E9").Value = "  +2.76%  "

$ws.Range("E10").Value = "  +0.64%  "

$ws.Range("D11").Value = "'0.0892"
$ws.Range("E11").Value = "  +0.73%  "

$ws.Range("D12").Value = "'1.877.99"
$ws.Range("E12").Value = "  +0.58%  "

$ws.Range("D13").Value = "'1.650.12"
$ws.Range("E13").Value = "  +0.82%  "

$ws.Range("D14").Value = "'4.06"
$ws.Range("E14").Value = "  +1.41%  "

$ws.Range("D15").Value = "'0.561"
$ws.Range("E15").Value = "  -3.16%  "

$ws.Range("D16").Value = "'64.74"
$ws.Range("E16").Value = "  +0.93%  "

$ws.Range("D17").Value = "'27.997.48"
$ws.Range("E17").Value = "  +2.12%  "

$ws.Range("D18").Value = "'233.73"
$ws.Range("E18").Value = "  +2.12%  "

$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "'7.66"
$ws.Range("E19").Value = "  +2.04%  "

$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "'0.0₃0724"
$ws.Range("E20").Value = "  +0.34%  "

$ws.Range("D21").Value = "'1.00"
$ws.Range("E21").Value = "  -0.20%  "

$ws.Range("E22").Value = "  +0.80%  "

$ws.Range("D23").Value = "'10.02"
$ws.Range("E23").Value = "  +3.87%  "

$ws.Range("E24").Value = "  +5.82%  "

$ws.Range("D25").Value = "'150.65"
$ws.Range("E25").Value = "  +0.98%  "

$ws.Range("E26").Value = "  -0.24%  "

$ws.Range("E27").Value = "  -0.58%  "

$ws.Range("D28").Value = "'15.71"
$ws.Range("E28").Value = "  +1.30%  "

$ws.Range("E29").Value = "  -0.23%  "

$ws.Range("E30").Value = "  +0.35%  "

$ws.Range("E31").Value = "  -0.61%  "

$ws.Range("E32").Value = "  +1.13%  "

$ws.Range("D33").Value = "'1.471.15"
$ws.Range("E33").Value = "  +4.44%  "

$ws.Range("E34").Value = "  -1.87%  "

$ws.Range("E35").Value = "  -1.82%  "

$ws.Range("E36").Value = "  -0.48%  "

$ws.Range("E37").Value = "  +0.10%  "

$ws.Range("D38").Value = "'0.883"
$ws.Range("E38").Value = "  +0.89%  "

$ws.Range("E39").Value = "  +0.96%  "

$ws.Range("E40").Value = "  +12.30%  "

$ws.Range("D41").Value = "'69.76"
$ws.Range("E41").Value = "  +7.86%  "

$ws.Range("E42").Value = "  -0.28%  "

$ws.Range("E43").Value = "  -2.03%  "

$ws.Range("E44").Value = "  -1.90%  "

$ws.Range("E45").Value = "  +0.48%  "

$ws.Range("E46").Value = "  -0.99%  "

$ws.Range("D47").Value = "'1.786.30"
$ws.Range("E47").Value = "  +0.52%  "

$ws.Range("E48").Value = "  +3.18%  "

$ws.Range("D49").Value = "'86.59"
$ws.Range("E49").Value = "  +0.95%  "

$ws.Range("E50").Value = "  -0.93%  "

$ws.Range("E51").Value = "  +0.42%  "
